$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style of an existing header cell (H1) onto the new headers so they
# match the bold/bordered/centered look of the other header cells
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ijData = @(
    @(6,6),
    @(8,8),
    @(5,6),
    @(7,7),
    @(8,8),
    @(9,9),
    @(6,6),
    @(7,7),
    @(6,6),
    @(7,7),
    @(7,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(6,6),
    @(7,7),
    @(7,7),
    @(7,7),
    @(9,9),
    @(8,8),
    @(7,7),
    @(5,6),
    @(7,7),
    @(9,9),
    @(9,9),
    @(8,8),
    @(4,5),
    @(7,7),
    @(5,6),
    @(5,6),
    @(7,7),
    @(3,4),
    @(6,6),
    @(5,5),
    @(3,4),
    @(8,8),
    @(6,7),
    @(7,7),
    @(6,6),
    @(3,4),
    @(4,4),
    @(7,7),
    @(7,7),
    @(7,8),
    @(6,6),
    @(7,7),
    @(8,8),
    @(7,8),
    @(7,7),
    @(7,7),
    @(6,6),
    @(7,7),
    @(9,9),
    @(8,8),
    @(9,9),
    @(8,8),
    @(6,7),
    @(7,7),
    @(8,8),
    @(7,8),
    @(6,6),
    @(6,6),
    @(6,6),
    @(8,8),
    @(7,8),
    @(8,8),
    @(6,7),
    @(8,9),
    @(7,7),
    @(9,9),
    @(8,8),
    @(9,9),
    @(6,6),
    @(7,8),
    @(5,5),
    @(8,8),
    @(4,5),
    @(9,9),
    @(9,9),
    @(6,6),
    @(6,6),
    @(4,4)
)

for ($idx = 0; $idx -lt $ijData.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $ijData[$idx][0]   # column I
    $ws.Cells.Item($row, 10).Value = $ijData[$idx][1]  # column J
}
